$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8, 11, 14: clear the "AL" column cell (discount amount no longer populated)
$ws.Range("AL8").Value = ""
$ws.Range("AL11").Value = ""
$ws.Range("AL14").Value = ""

# Row 15: quantity/amount figures corrected (order of magnitude fix),
# and the AL column value moves up into AK / AI while AL itself is cleared
$ws.Range("AH15").Value = 20
$ws.Range("AI15").Value = 300000
$ws.Range("AK15").Value = 300000
$ws.Range("AL15").Value = ""

# Row 16: same pattern as row 15
$ws.Range("AH16").Value = 500
$ws.Range("AI16").Value = 4000000
$ws.Range("AK16").Value = 4000000
$ws.Range("AL16").Value = ""
